$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.487.80'
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').Value = '3.985.60'
$ws.Range('E3').Value = '  -1.36%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '542.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.81'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.70%  '
$ws.Range('D7').Value = '3.976.64'
$ws.Range('E7').Value = '  -1.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.691'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  -3.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.168'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '56.80'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +19.36%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000320'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.79'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.51%  '
$ws.Range('D15').Value = '4.621.35'
$ws.Range('E15').Value = '  -1.47%  '
$ws.Range('D16').Value = '3.984.58'
$ws.Range('E16').Value = '  -1.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.01'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.61'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.66%  '
$ws.Range('E19').Value = '  -1.18%  '
$ws.Range('E20').Value = '  -2.79%  '
$ws.Range('D21').Value = '71.398.91'
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '429.65'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '98.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.23'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.56'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.51'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.76'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.77'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +16.73%  '
$ws.Range('E30').Value = '  +2.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '36.73'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.72'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +13.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '51.46'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +21.46%  '
$ws.Range('E34').Value = '  +1.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '13.44'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '684.70'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.77%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '65.66'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.442'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.152'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.74%  '
$ws.Range('E40').Value = '  -4.57%  '
$ws.Range('E41').Value = '  -3.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.29'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.84%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0487'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.16%  '
$ws.Range('E46').Value = '  -5.97%  '
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.74'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.52%  '
$ws.Range('E49').Value = '  -5.62%  '
$ws.Range('E50').Value = '  -1.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000273'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.96%  '
